# Update "想去人数" (want-to-go count) figures in column F for both the
# "展览" and "全部类型" sheets, which carry duplicate event listings.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    4  = 272
    8  = 2008
    10 = 98
    11 = 4317
    20 = 3081
    21 = 65
    26 = 77
    27 = 9
    30 = 198
    31 = 10
    32 = 485
    33 = 1702
    34 = 255
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
